$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.319.02"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "2.266.79"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.87"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.99"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.98%  "
$ws.Range("E7").Value = "  -2.29%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.539"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -5.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.38"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -6.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0829"
$ws.Range("D11").ClearFormats()
$ws.Range("E12").Value = "  -5.81%  "
$ws.Range("E13").Value = "  -2.47%  "
$ws.Range("D14").Value = "2.613.98"
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.850"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.87%  "
$ws.Range("D16").Value = "2.259.28"
$ws.Range("E16").Value = "  -0.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.07"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.55%  "
$ws.Range("D18").Value = "44.245.72"
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.21"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -5.03%  "
$ws.Range("D20").Value = "0.0₃0984"
$ws.Range("E20").Value = "  -1.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.40"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.87"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.39"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.02"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -6.13%  "
$ws.Range("E25").Value = "  -8.12%  "
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.06"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.65%  "
$ws.Range("E29").Value = "  -3.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.07"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -6.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.21"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.25"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.33%  "
$ws.Range("E33").Value = "  -4.37%  "
$ws.Range("E34").Value = "  +10.90%  "
$ws.Range("E35").Value = "  -3.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.111"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -5.35%  "
$ws.Range("E37").Value = "  -4.24%  "
$ws.Range("E38").Value = "  -2.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.54"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.54"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -10.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.97"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -10.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0310"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.80%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").Value = "1.738.83"
$ws.Range("E44").Value = "  -2.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "84.90"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.81%  "
$ws.Range("E46").Value = "  -6.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.23"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.80"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.56%  "
$ws.Range("E49").Value = "  -6.08%  "
$ws.Range("E50").Value = "  -5.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.45"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.80%  "
